# Burndown chart update - "grafico de burndown"
# Fill in actuals for weeks 4-6 (columns E, F, G) so the burndown chart's
# "Restante" series no longer shows #N/A for those points.

# NOTE: cells are filled one column at a time (all of column E, then all of
# column F, then all of column G) so the shared formulas in row 20 pick up
# each column's totals correctly before the next column is touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Column E (week 4): Requisito 2, 6, 8, 10, 11 each finish 1 point
$ws.Range("E3").Value = 1
$ws.Range("E7").Value = 1
$ws.Range("E9").Value = 1
$ws.Range("E11").Value = 1
$ws.Range("E12").Value = 1

# Column F (week 5): Requisito 2 finishes 5 points
$ws.Range("F3").Value = 5

# Column G (week 6): Requisito 6, 8, 10, 11 finish their remaining points
$ws.Range("G7").Value = 42
$ws.Range("G9").Value = 18
$ws.Range("G11").Value = 20
$ws.Range("G12").Value = 22

# Recalculate so the shared formulas in row 20 (Restante) resolve to numbers
# instead of #N/A now that SUM(E2:E19), SUM(F2:F19) and SUM(G2:G19) are > 0.
$excel.Calculate()

# Restore the view: scroll back to column A and move the selection to G13.
$ws.Activate()
$ws.Range("A1").Select()
$ws.Range("G13").Select()
